# Update automàtic: dades i banners [2026-02-24 18:20]
# Refresh DATA_EXTRACCIO timestamps and the handful of observation values
# that moved between the 17:48/17:50 and 18:18/18:20 meteo.cat scrapes.
# Percentage readings (column H) are written with a leading apostrophe so
# Excel keeps them as literal text ("68%") instead of coercing them into a
# numeric percentage value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-24 18:18:21"
$ws.Range("O2").Value = "6.1 °C"
$ws.Range("E3").Value = "2026-02-24 18:18:23"
$ws.Range("O3").Value = "4.3 °C"
$ws.Range("E4").Value = "2026-02-24 18:18:26"
$ws.Range("H4").Value = "'68%"
$ws.Range("E5").Value = "2026-02-24 18:18:28"
$ws.Range("E6").Value = "2026-02-24 18:18:31"
$ws.Range("E7").Value = "2026-02-24 18:18:33"
$ws.Range("J7").Value = "1020.4 hPa"
$ws.Range("E8").Value = "2026-02-24 18:18:36"
$ws.Range("J8").Value = "1019.8 hPa"
$ws.Range("E9").Value = "2026-02-24 18:18:38"
$ws.Range("E10").Value = "2026-02-24 18:18:41"
$ws.Range("E11").Value = "2026-02-24 18:18:43"
$ws.Range("O11").Value = "9.3 °C"
$ws.Range("E12").Value = "2026-02-24 18:18:45"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("E13").Value = "2026-02-24 18:18:47"
$ws.Range("J13").Value = "1023.6 hPa"
$ws.Range("O13").Value = "6.6 °C"
$ws.Range("E14").Value = "2026-02-24 18:18:50"
$ws.Range("H14").Value = "'82%"
$ws.Range("E15").Value = "2026-02-24 18:18:52"
$ws.Range("H15").Value = "'74%"
$ws.Range("E16").Value = "2026-02-24 18:18:54"
$ws.Range("O16").Value = "4.3 °C"
$ws.Range("E17").Value = "2026-02-24 18:18:57"
$ws.Range("H17").Value = "'34%"
$ws.Range("E18").Value = "2026-02-24 18:18:59"
$ws.Range("O18").Value = "11.2 °C"
$ws.Range("E19").Value = "2026-02-24 18:19:02"
$ws.Range("E20").Value = "2026-02-24 18:19:04"
$ws.Range("H20").Value = "'33%"
$ws.Range("E21").Value = "2026-02-24 18:19:07"
$ws.Range("H21").Value = "'60%"
$ws.Range("J21").Value = "1022.3 hPa"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-24 18:19:09"
$ws.Range("H22").Value = "'21%"
$ws.Range("E23").Value = "2026-02-24 18:19:12"
$ws.Range("E24").Value = "2026-02-24 18:19:14"
$ws.Range("J24").Value = "1021.6 hPa"
$ws.Range("O24").Value = "10.0 °C"
$ws.Range("E25").Value = "2026-02-24 18:19:17"
$ws.Range("E26").Value = "2026-02-24 18:19:19"
$ws.Range("O26").Value = "12.1 °C"
$ws.Range("E27").Value = "2026-02-24 18:19:22"
$ws.Range("H27").Value = "'31%"
$ws.Range("E28").Value = "2026-02-24 18:19:24"
$ws.Range("J28").Value = "1020.3 hPa"
$ws.Range("O28").Value = "11.7 °C"
$ws.Range("E29").Value = "2026-02-24 18:19:27"
$ws.Range("E30").Value = "2026-02-24 18:19:29"
$ws.Range("H30").Value = "'72%"
$ws.Range("J30").Value = "1020.0 hPa"
$ws.Range("E31").Value = "2026-02-24 18:19:31"
$ws.Range("H31").Value = "'59%"
$ws.Range("J31").Value = "1019.4 hPa"
$ws.Range("E32").Value = "2026-02-24 18:19:34"
$ws.Range("E33").Value = "2026-02-24 18:19:36"
$ws.Range("H33").Value = "'50%"
$ws.Range("J33").Value = "1021.9 hPa"
$ws.Range("O33").Value = "8.5 °C"
$ws.Range("E34").Value = "2026-02-24 18:19:39"
$ws.Range("O34").Value = "5.1 °C"
$ws.Range("E35").Value = "2026-02-24 18:19:41"
$ws.Range("H35").Value = "'38%"
$ws.Range("E36").Value = "2026-02-24 18:19:44"
$ws.Range("E37").Value = "2026-02-24 18:19:46"
$ws.Range("H37").Value = "'69%"
$ws.Range("E38").Value = "2026-02-24 18:19:49"
$ws.Range("E39").Value = "2026-02-24 18:19:51"
$ws.Range("H39").Value = "'36%"
$ws.Range("N39").Value = "1.6 °C 17:57 TU"
$ws.Range("E40").Value = "2026-02-24 18:19:54"
$ws.Range("J40").Value = "1022.9 hPa"
$ws.Range("O40").Value = "8.8 °C"
$ws.Range("E41").Value = "2026-02-24 18:19:56"
$ws.Range("H41").Value = "'76%"
$ws.Range("J41").Value = "1020.8 hPa"
$ws.Range("E42").Value = "2026-02-24 18:19:59"
$ws.Range("H42").Value = "'85%"
$ws.Range("E43").Value = "2026-02-24 18:20:01"
$ws.Range("O43").Value = "10.6 °C"
$ws.Range("E44").Value = "2026-02-24 18:20:03"
$ws.Range("H44").Value = "'38%"
$ws.Range("O44").Value = "2.9 °C"
$ws.Range("E45").Value = "2026-02-24 18:20:06"
$ws.Range("E46").Value = "2026-02-24 18:20:08"
$ws.Range("J46").Value = "1021.6 hPa"
$ws.Range("O46").Value = "10.3 °C"
